$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 167165
$ws.Range("I2").Value = 425
$ws.Range("J2").Value = 500645
$ws.Range("K2").Value = 425
$ws.Range("L2").Value = 500645
$ws.Range("M2").Value = -312
$ws.Range("N2").Value = -500871
$ws.Range("H4").Value = 72090.5
$ws.Range("I4").Value = 111240.89
$ws.Range("K4").Value = 111240.89
$ws.Range("M4").Value = -111126.89
$ws.Range("H6").Value = 772.8
$ws.Range("I6").Value = 10
$ws.Range("K6").Value = 30
$ws.Range("M6").Value = 82
$ws.Range("H9").Value = 90.818184
$ws.Range("I9").Value = 67.111115
$ws.Range("K9").Value = 67.111115
$ws.Range("M9").Value = 101.888885
$ws.Range("H43").Value = 97465.664
$ws.Range("I43").Value = 97465.664
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 97465.664
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -97396.664
$ws.Range("H51").Value = 5673.9165
$ws.Range("I51").Value = 5112.5713
$ws.Range("J51").Value = 6459.8
$ws.Range("K51").Value = 5112.5713
$ws.Range("L51").Value = 6459.8
$ws.Range("M51").Value = -4628.5713
$ws.Range("N51").Value = -7427.8
$ws.Range("H59").Value = 3583
$ws.Range("J59").Value = 4228.75
$ws.Range("L59").Value = 12686.25
$ws.Range("N59").Value = -13800.25
$ws.Range("H62").Value = 5385.5625
$ws.Range("I62").Value = 4764.3335
$ws.Range("J62").Value = 7249.25
$ws.Range("K62").Value = 4764.3335
$ws.Range("L62").Value = 7249.25
$ws.Range("M62").Value = -4140.3335
$ws.Range("N62").Value = -8497.25
$ws.Range("H65").Value = 5385.5625
$ws.Range("I65").Value = 4764.3335
$ws.Range("J65").Value = 7249.25
$ws.Range("K65").Value = 23821.6675
$ws.Range("L65").Value = 36246.25
$ws.Range("M65").Value = -20701.6675
$ws.Range("N65").Value = -42486.25
$ws.Range("H76").Value = 7784.3
$ws.Range("I76").Value = 7742.375
$ws.Range("K76").Value = 7742.375
$ws.Range("M76").Value = -7427.375
$ws.Range("H79").Value = 7784.3
$ws.Range("I79").Value = 7742.375
$ws.Range("K79").Value = 7742.375
$ws.Range("M79").Value = -6650.375
$ws.Range("I86").Value = 0
$ws.Range("J86").ClearContents()
$ws.Range("K86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -5246
$ws.Range("H88").Value = 15368.258
$ws.Range("J88").Value = 22843.25
$ws.Range("L88").Value = 22843.25
$ws.Range("N88").Value = -23655.25
$ws.Range("I89").Value = 0
$ws.Range("J89").ClearContents()
$ws.Range("K89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -26232
$ws.Range("H91").Value = 15368.258
$ws.Range("J91").Value = 22843.25
$ws.Range("L91").Value = 22843.25
$ws.Range("N91").Value = -25651.25
$ws.Range("H92").Value = 118214.5
$ws.Range("I92").Value = 50858
$ws.Range("K92").Value = 50858
$ws.Range("M92").Value = -49610
$ws.Range("H96").Value = 15351.895
$ws.Range("I96").Value = 29661.666
$ws.Range("K96").Value = 88984.99800000001
$ws.Range("M96").Value = -87611.99800000001
$ws.Range("H99").Value = 343.85715
$ws.Range("I99").Value = 343.85715
$ws.Range("K99").Value = 1031.57145
$ws.Range("M99").Value = 466.4285500000001
$ws.Range("H100").Value = 5741.6665
$ws.Range("I100").Value = 5834.5
$ws.Range("J100").Value = 4999
$ws.Range("K100").Value = 5834.5
$ws.Range("L100").Value = 4999
$ws.Range("M100").Value = -5293.5
$ws.Range("N100").Value = -6081
$ws.Range("H101").Value = 2661.4666
$ws.Range("I101").Value = 2195.5715
$ws.Range("J101").Value = 3069.125
$ws.Range("K101").Value = 6586.7145
$ws.Range("L101").Value = 9207.375
$ws.Range("M101").Value = -4964.7145
$ws.Range("N101").Value = -12451.375
$ws.Range("H106").Value = 4514.773
$ws.Range("I106").Value = 3772
$ws.Range("J106").Value = 5814.625
$ws.Range("K106").Value = 3772
$ws.Range("L106").Value = 5814.625
$ws.Range("M106").Value = -3141
$ws.Range("N106").Value = -7076.625
$ws.Range("H112").Value = 42938.863
$ws.Range("I112").Value = 1070.8572
$ws.Range("J112").Value = 62477.266
$ws.Range("K112").Value = 3212.5716
$ws.Range("L112").Value = 187431.798
$ws.Range("M112").Value = -2104.5716
$ws.Range("N112").Value = -189647.798
$ws.Range("H137").Value = 2333.4285
$ws.Range("I137").Value = 1607.3226
$ws.Range("K137").Value = 4821.9678
$ws.Range("M137").Value = -2271.9678
$ws.Range("H138").Value = 2667.0164
$ws.Range("I138").Value = 1149.2
$ws.Range("J138").Value = 3721.0557
$ws.Range("K138").Value = 3447.6
$ws.Range("L138").Value = 11163.1671
$ws.Range("M138").Value = 1692.4
$ws.Range("N138").Value = -21443.1671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 31333.447
$ws.Range("I32").Value = 34913.855
$ws.Range("J32").Value = 20890.584
$ws.Range("K32").Value = 34913.855
$ws.Range("L32").Value = 20890.584
$ws.Range("M32").Value = -34626.855
$ws.Range("N32").Value = -21464.584
$ws.Range("H97").Value = 874.7857
$ws.Range("I97").Value = 764
$ws.Range("J97").Value = 1281
$ws.Range("K97").Value = 764
$ws.Range("L97").Value = 1281
$ws.Range("M97").Value = -268
$ws.Range("N97").Value = -2273
$ws.Range("H122").Value = 50541.367
$ws.Range("I122").Value = 3253.1538
$ws.Range("J122").Value = 152999.17
$ws.Range("K122").Value = 9759.4614
$ws.Range("L122").Value = 458997.51
$ws.Range("M122").Value = -7309.4614
$ws.Range("N122").Value = -463897.51
$ws.Range("H133").Value = 65130.5
$ws.Range("J133").Value = 65130.5
$ws.Range("L133").Value = 65130.5
$ws.Range("N133").Value = -70190.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3236.7368
$ws.Range("I20").Value = 2289.5833
$ws.Range("J20").Value = 4860.4287
$ws.Range("K20").Value = 2289.5833
$ws.Range("L20").Value = 4860.4287
$ws.Range("M20").Value = -2042.5833
$ws.Range("N20").Value = -5354.4287
$ws.Range("H86").Value = 31641.768
$ws.Range("I86").Value = 2121.3845
$ws.Range("J86").Value = 54216.176
$ws.Range("K86").Value = 2121.3845
$ws.Range("L86").Value = 54216.176
$ws.Range("M86").Value = -998.3845000000001
$ws.Range("N86").Value = -56462.176
$ws.Range("H89").Value = 31641.768
$ws.Range("I89").Value = 2121.3845
$ws.Range("J89").Value = 54216.176
$ws.Range("K89").Value = 10606.9225
$ws.Range("L89").Value = 271080.88
$ws.Range("M89").Value = -4990.922500000001
$ws.Range("N89").Value = -282312.88
$ws.Range("H94").Value = 1089.2727
$ws.Range("I94").Value = 823.75
$ws.Range("K94").Value = 823.75
$ws.Range("M94").Value = -372.75
$ws.Range("H105").Value = 2994.1667
$ws.Range("I105").Value = 1998.5
$ws.Range("K105").Value = 1998.5
$ws.Range("M105").Value = -251.5
$ws.Range("H106").Value = 8281.429
$ws.Range("J106").Value = 8281.429
$ws.Range("L106").Value = 8281.429
$ws.Range("N106").Value = -10805.429

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1545.5385
$ws.Range("I16").Value = 1509.8
$ws.Range("K16").Value = 1509.8
$ws.Range("M16").Value = -1222.8
$ws.Range("H22").Value = 723.625
$ws.Range("J22").Value = 847.25
$ws.Range("L22").Value = 847.25
$ws.Range("N22").Value = -1547.25
$ws.Range("H58").Value = 7014.5415
$ws.Range("I58").Value = 5523.0586
$ws.Range("K58").Value = 5523.0586
$ws.Range("M58").Value = -5320.0586
$ws.Range("H62").Value = 68847.56
$ws.Range("I62").Value = 508502.5
$ws.Range("J62").Value = 6039.7144
$ws.Range("K62").Value = 508502.5
$ws.Range("L62").Value = 6039.7144
$ws.Range("M62").Value = -507878.5
$ws.Range("N62").Value = -7287.7144
$ws.Range("H65").Value = 68847.56
$ws.Range("I65").Value = 508502.5
$ws.Range("J65").Value = 6039.7144
$ws.Range("K65").Value = 2542512.5
$ws.Range("L65").Value = 30198.572
$ws.Range("M65").Value = -2539392.5
$ws.Range("N65").Value = -36438.572
$ws.Range("H93").Value = 22248.295
$ws.Range("I93").Value = 10656.454
$ws.Range("J93").Value = 43500
$ws.Range("K93").Value = 10656.454
$ws.Range("L93").Value = 43500
$ws.Range("M93").Value = -8784.454
$ws.Range("N93").Value = -47244
$ws.Range("H113").Value = 1545.5385
$ws.Range("I113").Value = 1509.8
$ws.Range("K113").Value = 1509.8
$ws.Range("M113").Value = 660.2
$ws.Range("H122").Value = 1998.7273
$ws.Range("I122").Value = 1757.8572
$ws.Range("J122").Value = 2420.25
$ws.Range("K122").Value = 5273.571599999999
$ws.Range("L122").Value = 7260.75
$ws.Range("M122").Value = -2823.571599999999
$ws.Range("N122").Value = -12160.75
$ws.Range("H124").Value = 42991.332
$ws.Range("J124").Value = 42991.332
$ws.Range("L124").Value = 42991.332
$ws.Range("N124").Value = -47901.332
$ws.Range("H136").Value = 7014.5415
$ws.Range("I136").Value = 5523.0586
$ws.Range("K136").Value = 16569.1758
$ws.Range("M136").Value = -14019.1758

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 624.5
$ws.Range("I2").Value = 97.5
$ws.Range("J2").Value = 1415
$ws.Range("K2").Value = 585
$ws.Range("L2").Value = 8490
$ws.Range("M2").Value = -472
$ws.Range("N2").Value = -8716
$ws.Range("H5").Value = 466.13043
$ws.Range("I5").Value = 354.52942
$ws.Range("J5").Value = 782.3333
$ws.Range("K5").Value = 1063.58826
$ws.Range("L5").Value = 2346.9999
$ws.Range("M5").Value = -951.58826
$ws.Range("N5").Value = -2570.9999
$ws.Range("H7").Value = 26
$ws.Range("I7").Value = 17.75
$ws.Range("J7").Value = 42.5
$ws.Range("K7").Value = 53.25
$ws.Range("L7").Value = 127.5
$ws.Range("M7").Value = 58.75
$ws.Range("N7").Value = -351.5
$ws.Range("H8").Value = 581.9231
$ws.Range("I8").Value = 581.9231
$ws.Range("K8").Value = 1745.7693
$ws.Range("M8").Value = -1606.7693
$ws.Range("H38").Value = 93.92856999999999
$ws.Range("I38").Value = 65
$ws.Range("J38").Value = 110
$ws.Range("K38").Value = 195
$ws.Range("L38").Value = 330
$ws.Range("M38").Value = 152
$ws.Range("N38").Value = -1024
$ws.Range("H68").Value = 17199.666
$ws.Range("I68").Value = 700
$ws.Range("J68").Value = 25449.5
$ws.Range("K68").Value = 2100
$ws.Range("L68").Value = 76348.5
$ws.Range("M68").Value = -1289
$ws.Range("N68").Value = -77970.5
$ws.Range("H71").Value = 17199.666
$ws.Range("I71").Value = 700
$ws.Range("J71").Value = 25449.5
$ws.Range("K71").Value = 6300
$ws.Range("L71").Value = 229045.5
$ws.Range("M71").Value = -2244
$ws.Range("N71").Value = -237157.5
$ws.Range("H92").Value = 664.2778
$ws.Range("I92").Value = 140.42857
$ws.Range("J92").Value = 997.63635
$ws.Range("K92").Value = 421.28571
$ws.Range("L92").Value = 2992.90905
$ws.Range("M92").Value = 826.71429
$ws.Range("N92").Value = -5488.90905
$ws.Range("H121").Value = 12766659
$ws.Range("I121").Value = 6391.2856
$ws.Range("J121").Value = 22691312
$ws.Range("K121").Value = 19173.8568
$ws.Range("L121").Value = 68073936
$ws.Range("M121").Value = -17863.8568
$ws.Range("N121").Value = -68076556
$ws.Range("H129").Value = 6498.8667
$ws.Range("I129").Value = 1099.8
$ws.Range("J129").Value = 9198.4
$ws.Range("K129").Value = 3299.4
$ws.Range("L129").Value = 27595.2
$ws.Range("M129").Value = 1700.6
$ws.Range("N129").Value = -37595.2
$ws.Range("H131").Value = 4957.923
$ws.Range("J131").Value = 6109.3335
$ws.Range("L131").Value = 18328.0005
$ws.Range("N131").Value = -28408.0005
$ws.Range("H135").Value = 466.13043
$ws.Range("I135").Value = 354.52942
$ws.Range("J135").Value = 782.3333
$ws.Range("K135").Value = 3190.76478
$ws.Range("L135").Value = 7040.9997
$ws.Range("M135").Value = -655.76478
$ws.Range("N135").Value = -12110.9997
$ws.Range("H136").Value = 2872.6428
$ws.Range("I136").Value = 2518.0833
$ws.Range("K136").Value = 7554.249899999999
$ws.Range("M136").Value = -2454.249899999999
$ws.Range("H141").Value = 3417.077
$ws.Range("I141").Value = 3417.077
$ws.Range("K141").Value = 10251.231
$ws.Range("M141").Value = -5071.231

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 15983
$ws.Range("I35").Value = 9999
$ws.Range("J35").Value = 18975
$ws.Range("K35").Value = 9999
$ws.Range("L35").Value = 18975
$ws.Range("M35").Value = -9701
$ws.Range("N35").Value = -19571
$ws.Range("H36").Value = 72997.664
$ws.Range("I36").Value = 115196
$ws.Range("J36").Value = 20249.75
$ws.Range("K36").Value = 115196
$ws.Range("L36").Value = 20249.75
$ws.Range("M36").Value = -114711
$ws.Range("N36").Value = -21219.75
$ws.Range("H39").Value = 100555.5
$ws.Range("J39").Value = 100555.5
$ws.Range("L39").Value = 100555.5
$ws.Range("N39").Value = -101619.5
$ws.Range("H43").Value = 18882.637
$ws.Range("I43").Value = 6927.5
$ws.Range("J43").Value = 25714.143
$ws.Range("K43").Value = 6927.5
$ws.Range("L43").Value = 25714.143
$ws.Range("M43").Value = -6776.5
$ws.Range("N43").Value = -26016.143
$ws.Range("H80").Value = 21543.715
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 21543.715
$ws.Range("K80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("M80").Value = 21543.715
$ws.Range("N80").Value = -23539.715
$ws.Range("H83").Value = 21543.715
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 21543.715
$ws.Range("K83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("M83").Value = 107718.575
$ws.Range("N83").Value = -117702.575
$ws.Range("H93").Value = 80000
$ws.Range("J93").Value = 80000
$ws.Range("L93").Value = 80000
$ws.Range("N93").Value = -83744
$ws.Range("H97").Value = 1329.4445
$ws.Range("I97").Value = 863.1667
$ws.Range("J97").Value = 2262
$ws.Range("K97").Value = 863.1667
$ws.Range("L97").Value = 2262
$ws.Range("M97").Value = -367.1667
$ws.Range("N97").Value = -3254
$ws.Range("H102").Value = 36110.31
$ws.Range("I102").Value = 42980.5
$ws.Range("K102").Value = 42980.5
$ws.Range("M102").Value = -41358.5
$ws.Range("H106").Value = 41139.5
$ws.Range("J106").Value = 41139.5
$ws.Range("L106").Value = 41139.5
$ws.Range("N106").Value = -43663.5
$ws.Range("H123").Value = 86806.14
$ws.Range("J123").Value = 86806.14
$ws.Range("L123").Value = 86806.14
$ws.Range("N123").Value = -91706.14
$ws.Range("H136").Value = 15121.423
$ws.Range("J136").Value = 15121.423
$ws.Range("L136").Value = 45364.269
$ws.Range("N136").Value = -50464.269
$ws.Range("H140").Value = 137139.6
$ws.Range("J140").Value = 137139.6
$ws.Range("L140").Value = 137139.6
$ws.Range("N140").Value = -147499.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2018.3226
$ws.Range("I16").Value = 1273.875
$ws.Range("J16").Value = 4570.7144
$ws.Range("K16").Value = 1273.875
$ws.Range("L16").Value = 4570.7144
$ws.Range("M16").Value = -1103.875
$ws.Range("N16").Value = -4910.7144
$ws.Range("H22").Value = 7442.9624
$ws.Range("I22").Value = 2761.5
$ws.Range("J22").Value = 7825.1226
$ws.Range("K22").Value = 2761.5
$ws.Range("L22").Value = 7825.1226
$ws.Range("M22").Value = -2466.5
$ws.Range("N22").Value = -8415.122599999999
$ws.Range("H27").Value = 7442.9624
$ws.Range("I27").Value = 2761.5
$ws.Range("J27").Value = 7825.1226
$ws.Range("K27").Value = 2761.5
$ws.Range("L27").Value = 7825.1226
$ws.Range("M27").Value = -2654.5
$ws.Range("N27").Value = -8039.1226
$ws.Range("H61").Value = 8067.5
$ws.Range("I61").Value = 7641
$ws.Range("J61").Value = 10200
$ws.Range("K61").Value = 7641
$ws.Range("L61").Value = 10200
$ws.Range("M61").Value = -7439
$ws.Range("N61").Value = -10604
$ws.Range("H74").Value = 68287.57000000001
$ws.Range("I74").Value = 59759.4
$ws.Range("K74").Value = 59759.4
$ws.Range("M74").Value = -58761.4
$ws.Range("H77").Value = 68287.57000000001
$ws.Range("I77").Value = 59759.4
$ws.Range("K77").Value = 179278.2
$ws.Range("M77").Value = -174286.2
$ws.Range("H93").Value = 3056.875
$ws.Range("I93").Value = 2735.6
$ws.Range("J93").Value = 3592.3333
$ws.Range("K93").Value = 2735.6
$ws.Range("L93").Value = 3592.3333
$ws.Range("M93").Value = -1487.6
$ws.Range("N93").Value = -6088.3333
$ws.Range("H98").Value = 49071
$ws.Range("J98").Value = 49071
$ws.Range("L98").Value = 49071
$ws.Range("N98").Value = -55061
$ws.Range("H113").Value = 8067.5
$ws.Range("I113").Value = 7641
$ws.Range("J113").Value = 10200
$ws.Range("K113").Value = 7641
$ws.Range("L113").Value = 10200
$ws.Range("M113").Value = -5471
$ws.Range("N113").Value = -14540
$ws.Range("H122").Value = 3820.0356
$ws.Range("I122").Value = 3596.423
$ws.Range("J122").Value = 6727
$ws.Range("K122").Value = 10789.269
$ws.Range("L122").Value = 20181
$ws.Range("M122").Value = -8339.269
$ws.Range("N122").Value = -25081

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 65599.39999999999
$ws.Range("I64").Value = 49999
$ws.Range("K64").Value = 49999
$ws.Range("M64").Value = -49751
$ws.Range("H67").Value = 65599.39999999999
$ws.Range("I67").Value = 49999
$ws.Range("K67").Value = 49999
$ws.Range("M67").Value = -49141
$ws.Range("H96").Value = 3694.9333
$ws.Range("I96").Value = 4430.3335
$ws.Range("K96").Value = 4430.3335
$ws.Range("M96").Value = -3057.3335
$ws.Range("H97").Value = 24151.834
$ws.Range("J97").Value = 24151.834
$ws.Range("L97").Value = 24151.834
$ws.Range("N97").Value = -26133.834
$ws.Range("H100").Value = 1084.3
$ws.Range("I100").Value = 1018.9167
$ws.Range("J100").Value = 1345.8334
$ws.Range("K100").Value = 2037.8334
$ws.Range("L100").Value = 2691.6668
$ws.Range("M100").Value = -1496.8334
$ws.Range("N100").Value = -3773.6668
$ws.Range("H104").Value = 52777.777
$ws.Range("J104").Value = 52777.777
$ws.Range("L104").Value = 52777.777
$ws.Range("N104").Value = -59765.777
$ws.Range("H132").Value = 4684.1665
$ws.Range("I132").Value = 3644.9575
$ws.Range("J132").Value = 8441.308000000001
$ws.Range("K132").Value = 10934.8725
$ws.Range("L132").Value = 25323.924
$ws.Range("M132").Value = -8404.872499999999
$ws.Range("N132").Value = -30383.924
$ws.Range("H136").Value = 8552
$ws.Range("I136").Value = 8830.333000000001
$ws.Range("K136").Value = 26490.999
$ws.Range("M136").Value = -23940.999
